# Apply Coin/Link/Price/Volume(1h) refresh from GitHub Actions scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.123.97"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "3.581.65"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'574.67"
$ws.Range("E5").Value = "  -5.48%  "
$ws.Range("D6").Value = "'192.16"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "3.574.40"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("D8").Value = "'0.616"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.678"
$ws.Range("E10").Value = "  -5.72%  "
$ws.Range("D11").Value = "'55.83"
$ws.Range("E11").Value = "  -6.30%  "
$ws.Range("E12").Value = "  -6.28%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("D14").Value = "'9.85"
$ws.Range("E14").Value = "  -4.95%  "
$ws.Range("D15").Value = "4.157.87"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").Value = "3.584.23"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "'0.125"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "'18.37"
$ws.Range("E18").Value = "  -4.81%  "
$ws.Range("D19").Value = "67.084.69"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").Value = "'12.17"
$ws.Range("E20").Value = "  -4.65%  "
$ws.Range("E21").Value = "  -6.64%  "
$ws.Range("D22").Value = "'400.16"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("E23").Value = "  -8.07%  "
$ws.Range("D24").Value = "'85.97"
$ws.Range("E24").Value = "  -4.20%  "
$ws.Range("D25").Value = "'11.39"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "'2.93"
$ws.Range("E26").Value = "  -3.78%  "
$ws.Range("E27").Value = "  -3.77%  "
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "'3.60"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "'8.93"
$ws.Range("E30").Value = "  -6.69%  "
$ws.Range("D31").Value = "'7.66"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'31.21"
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("D33").Value = "'637.93"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "'12.13"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("E35").Value = "  -5.67%  "
$ws.Range("D36").Value = "'64.01"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("D37").Value = "'42.26"
$ws.Range("E37").Value = "  -10.97%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "0.0₃0772"
$ws.Range("E40").Value = "  -6.34%  "
$ws.Range("D41").Value = "3.192.87"
$ws.Range("E41").Value = "  +10.69%  "
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'2.97"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.69"
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("D46").Value = "'0.0416"
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.10"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.130"
$ws.Range("E48").Value = "  -6.33%  "
$ws.Range("D49").Value = "'142.69"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "'8.59"
$ws.Range("E51").Value = "  -6.41%  "
